# FINAL DEL 18 SEPT 2021
$wb = $excel.ActiveWorkbook

# The workbook opens with the "ARQUITECTO" sheet (tabSelected="1") active;
# that is the sheet whose vale amount/selection changed.
$ws = $wb.ActiveSheet

# Written-out amount: "CIENTO CINCUENTA" (150) -> "CIEn" (100).
$ws.Range("A2").Value = "CIEn    MIL   PESOS 00/100 M.N."

# Numeric amount backing the vale: 150000 -> 100000.
$ws.Range("D1").Value = 100000

# Move the saved cursor/selection to C12.
[void]$ws.Range("C12").Select()

# Recalculate so the volatile TODAY() cells (row 11, both sheets) pick up
# the current "as of" date used when this vale was finalized.
[void]$excel.Calculate()
